$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.126.72'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.812.97'
$ws.Range("E3").Value = '  -1.74%  '
$ws.Range("E4").Value = '  +0.71%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.16'
$ws.Range("E5").Value = '  +2.18%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  +0.81%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.58'
$ws.Range("E8").Value = '  -6.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.325'
$ws.Range("E9").Value = '  +5.74%  '
$ws.Range("E10").Value = '  -1.36%  '
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.076.52'
$ws.Range("E12").Value = '  -1.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.809.16'
$ws.Range("E13").Value = '  -2.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.666'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '11.06'
$ws.Range("E15").Value = '  -4.98%  '
$ws.Range("E16").Value = '  -2.51%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.083.24'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.69'
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0792'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '238.61'
$ws.Range("E20").Value = '  -3.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.97'
$ws.Range("E21").Value = '  -2.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.69'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.26'
$ws.Range("E24").Value = '  +2.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.86'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.85'
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.54'
$ws.Range("E27").Value = '  -2.27%  '
$ws.Range("E28").Value = '  -1.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.60'
$ws.Range("E29").Value = '  +20.97%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.16'
$ws.Range("E31").Value = '  +5.48%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0553'
$ws.Range("E32").Value = '  +2.63%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.03'
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.78'
$ws.Range("E34").Value = '  -5.71%  '
$ws.Range("B35").Value = 'TrustWalletToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.14'
$ws.Range("E35").Value = '  +5.20%  '
$ws.Range("B36").Value = 'Aave'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '92.41'
$ws.Range("E36").Value = '  +2.06%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.679'
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0194'
$ws.Range("E38").Value = '  -0.57%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.313.53'
$ws.Range("E39").Value = '  -2.26%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -2.55%  '
$ws.Range("E41").Value = '  +1.44%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.59'
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.47'
$ws.Range("E43").Value = '  +1.24%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.30'
$ws.Range("E44").Value = '  -5.65%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.76'
$ws.Range("E45").Value = '  -2.23%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.31'
$ws.Range("E46").Value = '  +3.73%  '
$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0511'
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.991.79'
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.01'
$ws.Range("E49").Value = '  +0.69%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0653'
$ws.Range("E50").Value = '  +4.84%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '99.28'
$ws.Range("E51").Value = '  -5.34%  '
